# Swap the "JUMLAH" (A) and "TAHUN" (B) columns so that TAHUN (year) is in
# column A and JUMLAH (count) is in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns A ("JUMLAH") and B ("TAHUN") so TAHUN (year) becomes the
# first column and JUMLAH (count) becomes the second one. Use a real
# cut + insert (column-wide), the same gesture a user performs in the UI
# ("select column B, Cut, select column A, Insert Cut Cells") - this keeps
# each row's existing cell style attached to the row instead of resetting
# it, unlike a plain value copy/paste.
$ws.Columns.Item(2).Cut()
$ws.Columns.Item(1).Insert()

# Update the active cell selection.
$ws.Range("E9").Select()
